# Fruta / hortaliza, semanal
# Insert this week's two new price records (rows 51-52) for Frambuesa at
# Vega Central Mapocho de Santiago, pushing the existing history down by
# two rows (old row 51 -> new row 53, ... old row 72 -> new row 74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (originally rows 51-72) down by two rows so the
# two new weekly records can be placed at the top of the block.
$ws.Rows("51:52").Insert()

# --- New row 51 ------------------------------------------------------
$ws.Range("A51").Value = 9
$ws.Range("B51").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C51").Value = "Metropolitana"
$ws.Range("D51").Value = 44582
$ws.Range("E51").Value = 13
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100101
$ws.Range("H51").Value = "Berries"
$ws.Range("I51").Value = 100101004
$ws.Range("J51").Value = "Frambuesa"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Especial"
$ws.Range("M51").Value = 440
$ws.Range("N51").Value = 8000
$ws.Range("O51").Value = 8000
$ws.Range("P51").Value = 8000
$ws.Range("Q51").Value = "$/bandeja 2 kilos"
$ws.Range("R51").Value = "Provincia de Linares"
$ws.Range("S51").Value = 4000
$ws.Range("T51").Value = 2

# --- New row 52 ------------------------------------------------------
$ws.Range("A52").Value = 9
$ws.Range("B52").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 44582
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100101
$ws.Range("H52").Value = "Berries"
$ws.Range("I52").Value = 100101004
$ws.Range("J52").Value = "Frambuesa"
$ws.Range("K52").Value = "Sin especificar"
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 380
$ws.Range("N52").Value = 7000
$ws.Range("O52").Value = 7000
$ws.Range("P52").Value = 7000
$ws.Range("Q52").Value = "$/bandeja 2 kilos"
$ws.Range("R52").Value = "Provincia de Linares"
$ws.Range("S52").Value = 3500
$ws.Range("T52").Value = 2
